$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row -> Column -> new value
$updates = @{
    2  = @{ C = 347; F = 45597; G = 30468; H = 45658 }
    3  = @{ E = 30407; F = 45597 }
    4  = @{ E = 30407; F = 45597; G = 30468; H = 45658 }
    5  = @{ E = 30407; F = 45597; G = 30468; H = 45658 }
    6  = @{ C = 455; F = 45566; G = 30468; H = 45658 }
    7  = @{ E = 30376; F = 45566; G = 30468; H = 45658 }
    8  = @{ D = 423; E = 30376; F = 45566; H = 45658 }
    9  = @{ E = 30407; F = 45597; G = 30468; H = 45658 }
    10 = @{ D = 496; E = 30407; F = 45597; H = 45658 }
    11 = @{ E = 30376; F = 45566; G = 30468; H = 45658 }
    12 = @{ C = 384; D = 365; F = 45597; H = 45658 }
    13 = @{ C = 480; F = 45597; G = 30468; H = 45658 }
    14 = @{ C = 435; D = 409; F = 45566; H = 45658 }
    15 = @{ G = 30438; H = 45658 }
    16 = @{ C = 467; D = 423; F = 45566; H = 45658 }
    17 = @{ C = 383; D = 407; F = 45597; H = 45658 }
    18 = @{ D = 273; E = 30407; F = 45597; H = 45658 }
    19 = @{ D = 411; E = 30407; F = 45597; H = 45658 }
    20 = @{ C = 493; F = 45566; G = 30468; H = 45658 }
    21 = @{ C = 336; D = 380; F = 45597; H = 45658 }
    22 = @{ G = 30468; H = 45658 }
    23 = @{ D = 393; E = 30376; F = 45566; H = 45658 }
    24 = @{ D = 230; E = 30407; F = 45597; H = 45658 }
    25 = @{ C = 406; F = 45566; G = 30468; H = 45658 }
    26 = @{ E = 30407; F = 45597; G = 30468; H = 45658 }
    27 = @{ D = 411; H = 45658 }
    28 = @{ D = 423; E = 30407; F = 45597; H = 45658 }
    29 = @{ C = 372; D = 380; F = 45597; H = 45658 }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $cellRef = "$colLetter$rowNum"
        $ws.Range($cellRef).Value = $cols[$colLetter]
    }
}
